$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (GitHub Actions scheduled refresh)

# Row 2
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.708.65'
$ws.Range('D2').ClearFormats()

# Row 3
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.600.18'
$ws.Range('D3').ClearFormats()

# Row 4
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').ClearFormats()

# Row 5
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.27'
$ws.Range('D5').ClearFormats()

# Row 6
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.511'
$ws.Range('D6').ClearFormats()

# Row 7
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('D7').ClearFormats()

# Row 8
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0619'
$ws.Range('D8').ClearFormats()

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.247'
$ws.Range('D9').ClearFormats()

# Row 10
$ws.Range('E10').Value = '  +1.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.70'
$ws.Range('D10').ClearFormats()

# Row 11
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0840'
$ws.Range('D11').ClearFormats()

# Row 12
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.824.83'
$ws.Range('D12').ClearFormats()

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.611.08'
$ws.Range('D13').ClearFormats()

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('D14').ClearFormats()

# Row 15
$ws.Range('E15').Value = '  -0.04%  '

# Row 16
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.16'
$ws.Range('D16').ClearFormats()

# Row 17
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.672.22'
$ws.Range('D17').ClearFormats()

# Row 18
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0746'
$ws.Range('D18').ClearFormats()

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.22'
$ws.Range('D19').ClearFormats()

# Row 20
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '209.84'
$ws.Range('D20').ClearFormats()

# Row 21
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').ClearFormats()

# Row 22
$ws.Range('E22').Value = '  +0.79%  '

# Row 23
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.29'
$ws.Range('D23').ClearFormats()

# Row 24
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.95'
$ws.Range('D24').ClearFormats()

# Row 25
$ws.Range('E25').Value = '  -1.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.87'
$ws.Range('D25').ClearFormats()

# Row 26
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()

# Row 27
$ws.Range('E27').Value = '  -0.94%  '

# Row 28
$ws.Range('E28').Value = '  -0.90%  '

# Row 29
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.44'
$ws.Range('D29').ClearFormats()

# Row 30
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0514'
$ws.Range('D30').ClearFormats()

# Row 31
$ws.Range('E31').Value = '  -0.68%  '

# Row 32
$ws.Range('E32').Value = '  +0.81%  '

# Row 33
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('D33').ClearFormats()

# Row 34
$ws.Range('E34').Value = '  +0.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.294.77'
$ws.Range('D34').ClearFormats()

# Row 35
$ws.Range('E35').Value = '  +0.80%  '

# Row 36
$ws.Range('E36').Value = '  -4.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.606'
$ws.Range('D36').ClearFormats()

# Row 37
$ws.Range('E37').Value = '  +0.31%  '

# Row 38
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0170'
$ws.Range('D38').ClearFormats()

# Row 39
$ws.Range('E39').Value = '  +17.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.10'
$ws.Range('D39').ClearFormats()

# Row 40
$ws.Range('E40').Value = '  -2.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.822'
$ws.Range('D40').ClearFormats()

# Row 41
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.23'
$ws.Range('D41').ClearFormats()

# Row 42
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E42').Value = '  -1.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.41'
$ws.Range('D42').ClearFormats()

# Row 43
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.780'
$ws.Range('D43').ClearFormats()

# Row 44
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.13'
$ws.Range('D44').ClearFormats()

# Row 45
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.735.85'
$ws.Range('D45').ClearFormats()

# Row 46
$ws.Range('E46').Value = '  +1.00%  '

# Row 47
$ws.Range('E47').Value = '  -2.82%  '

# Row 48
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('D48').ClearFormats()

# Row 49
$ws.Range('E49').Value = '  -1.73%  '

# Row 50
$ws.Range('E50').Value = '  +1.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0515'
$ws.Range('D50').ClearFormats()

# Row 51
$ws.Range('E51').Value = '  +0.20%  '
